$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.606343030929565
$ws.Range("B1").Value = 1.391809821128845
$ws.Range("C1").Value = 5.278607845306396
$ws.Range("D1").Value = 3.333710670471191
$ws.Range("E1").Value = 0.5820793509483337
